$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TMALocations")

# Insert a new row at position 65 (everything below shifts down by one)
$null = $ws.Rows.Item(65).Insert()

# Populate the new row: branch 317 "HILLSBORO HOOP SHED", region CGS
$ws.Cells.Item(65, 1).Value = 317
$ws.Cells.Item(65, 2).Value = "HILLSBORO HOOP SHED"
$ws.Cells.Item(65, 4).Value = "CGS"

# Grow the TMALocations table to include the newly inserted row
$lo = $ws.ListObjects.Item(1)
$null = $lo.Resize($ws.Range("A1:M105"))

# Make TMALocations the active sheet/tab and select the new row's F cell
$null = $ws.Activate()
$null = $ws.Range("F65").Select()
